# "Finish add education and add certifications in the profile page"
#
# Semantic changes applied:
#  1. AddEduction!A2 text "NZ" -> "New Zealand" (shared-string table gains a
#     new entry; the other shared cells that merely shift index because the
#     old "NZ" entry was removed from the middle of the table come along for
#     free once the host re-serialises sharedStrings.xml).
#  2. Selection/active-cell bookkeeping:
#       - AddNewLanguages: no longer the tab shown on open; selection moves
#         off B9 -> C33.
#       - AddEduction: selection moves off E15 -> A6.
#       - AddCertifications: becomes the tab shown on open (tabSelected);
#         selection moves off F14 -> A2.
#  3. AddCertifications column layout: column A narrows and a width is now
#     recorded for column D (room made for the "certification" detail work).

$wb = $excel.ActiveWorkbook

$wsLang = $wb.Worksheets.Item("AddNewLanguages")
$wsEdu  = $wb.Worksheets.Item("AddEduction")
$wsCert = $wb.Worksheets.Item("AddCertifications")

# --- 1. Content fix: "NZ" -> "New Zealand" on the AddEduction sheet -------
$wsEdu.Range("A2").Value = "New Zealand"

# --- 2. Column width tweaks on AddCertifications ---------------------------
# The host's ColumnWidth setter quantises to an MDW-7 pixel grid, so these
# land on the closest achievable grid point to the authored
# 18.73046875 / 24.9296875 raw widths.
$wsCert.Columns.Item(1).ColumnWidth = 18
$wsCert.Columns.Item(4).ColumnWidth = 24.285714285714285

# --- 3. Selections + which tab is active when the workbook is reopened -----
# Visit AddNewLanguages and AddEduction first (leaving behind their new
# selections, and no longer the active/tabSelected sheet), then finish on
# AddCertifications so it ends up both ActiveSheet and tabSelected.
$wsLang.Activate()
$wsLang.Range("C33").Select()

$wsEdu.Activate()
$wsEdu.Range("A6").Select()

$wsCert.Activate()
$wsCert.Range("A2").Select()
